# Apply scheduled-runner profit/price updates to the Ultima Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 792.8214
$ws.Range("I80").Value = 1453.3334
$ws.Range("J80").Value = 479.94736
$ws.Range("K80").Value = 4360.0002
$ws.Range("L80").Value = 1439.84208
$ws.Range("M80").Value = -3362.0002
$ws.Range("N80").Value = -3435.84208

$ws.Range("H83").Value = 792.8214
$ws.Range("I83").Value = 1453.3334
$ws.Range("J83").Value = 479.94736
$ws.Range("K83").Value = 13080.0006
$ws.Range("L83").Value = 4319.52624
$ws.Range("M83").Value = -8088.000599999999
$ws.Range("N83").Value = -14303.52624

$ws.Range("H138").Value = 6946281.5
$ws.Range("I138").Value = 9525426
$ws.Range("J138").Value = 2430
$ws.Range("K138").Value = 28576278
$ws.Range("L138").Value = 7290
$ws.Range("M138").Value = -28571138
$ws.Range("N138").Value = -17570

$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2334.1333
$ws.Range("I45").Value = 2271.2
$ws.Range("J45").Value = 2460
$ws.Range("K45").Value = 2271.2
$ws.Range("L45").Value = 2460
$ws.Range("M45").Value = -1894.2
$ws.Range("N45").Value = -3214

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1346.6666
$ws.Range("I107").Value = 1107.625
$ws.Range("J107").Value = 1824.75
$ws.Range("K107").Value = 1107.625
$ws.Range("L107").Value = 1824.75
$ws.Range("M107").Value = 812.375
$ws.Range("N107").Value = -5664.75

$ws.Range("H134").Value = 2887.6047
$ws.Range("I134").Value = 1869.7407
$ws.Range("J134").Value = 4605.25
$ws.Range("K134").Value = 5609.2221
$ws.Range("L134").Value = 13815.75
$ws.Range("M134").Value = -3074.2221
$ws.Range("N134").Value = -18885.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 100.454544
$ws.Range("I7").Value = 71.666664
$ws.Range("J7").Value = 135
$ws.Range("K7").Value = 71.666664
$ws.Range("L7").Value = 135
$ws.Range("M7").Value = 41.333336
$ws.Range("N7").Value = -361

$ws.Range("H20").Value = 35450
$ws.Range("J20").Value = 35450
$ws.Range("L20").Value = 35450
$ws.Range("N20").Value = -35922

$ws.Range("H22").Value = 516.25
$ws.Range("I22").Value = 361.875
$ws.Range("J22").Value = 825
$ws.Range("K22").Value = 361.875
$ws.Range("L22").Value = 825
$ws.Range("M22").Value = -11.875
$ws.Range("N22").Value = -1525

$ws.Range("H30").Value = 35450
$ws.Range("J30").Value = 35450
$ws.Range("L30").Value = 35450
$ws.Range("N30").Value = -35632

$ws.Range("H31").Value = 8777629
$ws.Range("I31").Value = 7177.2173
$ws.Range("J31").Value = 22225654
$ws.Range("K31").Value = 7177.2173
$ws.Range("L31").Value = 22225654
$ws.Range("M31").Value = -6882.2173
$ws.Range("N31").Value = -22226244

$ws.Range("H34").Value = 8777629
$ws.Range("I34").Value = 7177.2173
$ws.Range("J34").Value = 22225654
$ws.Range("K34").Value = 7177.2173
$ws.Range("L34").Value = 22225654
$ws.Range("M34").Value = -6975.2173
$ws.Range("N34").Value = -22226058

$ws.Range("H64").Value = 23270
$ws.Range("J64").Value = 23270
$ws.Range("L64").Value = 23270
$ws.Range("N64").Value = -23766

$ws.Range("H67").Value = 23270
$ws.Range("J67").Value = 23270
$ws.Range("L67").Value = 23270
$ws.Range("N67").Value = -24986

$ws.Range("H86").Value = 2866.2632
$ws.Range("I86").Value = 2515.3845
$ws.Range("K86").Value = 2515.3845
$ws.Range("M86").Value = -1392.3845

$ws.Range("H89").Value = 2866.2632
$ws.Range("I89").Value = 2515.3845
$ws.Range("K89").Value = 12576.9225
$ws.Range("M89").Value = -6960.922500000001

$ws.Range("H107").Value = 1122
$ws.Range("I107").Value = 1202.5
$ws.Range("J107").Value = 800
$ws.Range("K107").Value = 1202.5
$ws.Range("L107").Value = 800
$ws.Range("M107").Value = 717.5
$ws.Range("N107").Value = -4640

$ws.Range("H128").Value = 35450
$ws.Range("J128").Value = 35450
$ws.Range("L128").Value = 35450
$ws.Range("N128").Value = -45410

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 825.0516
$ws.Range("I131").Value = 483.33334
$ws.Range("J131").Value = 860
$ws.Range("K131").Value = 1450.00002
$ws.Range("L131").Value = 2580
$ws.Range("M131").Value = 3589.99998
$ws.Range("N131").Value = -12660

$ws.Range("H132").Value = 1441.9166
$ws.Range("I132").Value = 436
$ws.Range("J132").Value = 1944.875
$ws.Range("K132").Value = 3924
$ws.Range("L132").Value = 17503.875
$ws.Range("M132").Value = -1394
$ws.Range("N132").Value = -22563.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3142.4849
$ws.Range("I102").Value = 3336.3572
$ws.Range("J102").Value = 2056.8
$ws.Range("K102").Value = 3336.3572
$ws.Range("L102").Value = 2056.8
$ws.Range("M102").Value = -1714.3572
$ws.Range("N102").Value = -5300.8

$ws.Range("H126").Value = 3922.84
$ws.Range("I126").Value = 2131.0908
$ws.Range("J126").Value = 5330.643
$ws.Range("K126").Value = 6393.2724
$ws.Range("L126").Value = 15991.929
$ws.Range("M126").Value = -3923.2724
$ws.Range("N126").Value = -20931.929

$ws.Range("H132").Value = 5062.5654
$ws.Range("I132").Value = 4431.2334
$ws.Range("J132").Value = 6246.3125
$ws.Range("K132").Value = 13293.7002
$ws.Range("L132").Value = 18738.9375
$ws.Range("M132").Value = -10763.7002
$ws.Range("N132").Value = -23798.9375

$ws.Range("H136").Value = 24330.4
$ws.Range("J136").Value = 24330.4
$ws.Range("L136").Value = 72991.20000000001
$ws.Range("N136").Value = -78091.20000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5291.7046
$ws.Range("I7").Value = 5406.0415
$ws.Range("J7").Value = 5154.5
$ws.Range("K7").Value = 5406.0415
$ws.Range("L7").Value = 5154.5
$ws.Range("M7").Value = -5294.0415
$ws.Range("N7").Value = -5378.5

$ws.Range("H40").Value = 6233.294
$ws.Range("I40").Value = 8419.75
$ws.Range("J40").Value = 4289.778
$ws.Range("K40").Value = 8419.75
$ws.Range("L40").Value = 4289.778
$ws.Range("M40").Value = -8283.75
$ws.Range("N40").Value = -4561.778

$ws.Range("H126").Value = 5291.7046
$ws.Range("I126").Value = 5406.0415
$ws.Range("J126").Value = 5154.5
$ws.Range("K126").Value = 16218.1245
$ws.Range("L126").Value = 15463.5
$ws.Range("M126").Value = -13748.1245
$ws.Range("N126").Value = -20403.5

$ws.Range("H139").Value = 56216.668
$ws.Range("J139").Value = 56216.668
$ws.Range("L139").Value = 56216.668
$ws.Range("N139").Value = -66496.66800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3923.0625
$ws.Range("I126").Value = 2828
$ws.Range("J126").Value = 8668.333000000001
$ws.Range("K126").Value = 8484
$ws.Range("L126").Value = 26004.999
$ws.Range("M126").Value = -6014
$ws.Range("N126").Value = -30944.999

$ws.Range("H136").Value = 1036.7297
$ws.Range("I136").Value = 774.82355
$ws.Range("K136").Value = 2324.47065
$ws.Range("M136").Value = 225.5293500000002
